$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns I ("I0") and J ("IF"), matching the style
# of the existing header cells (B1:H1 use a bold/bordered/centered format).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$values = @{
    2  = @(5, 6)
    3  = @(7, 7)
    4  = @(4, 6)
    5  = @(1, 3)
    6  = @(1, 4)
    7  = @(1, 3)
    8  = @(1, 4)
    9  = @(1, 2)
    10 = @(7, 9)
    11 = @(7, 7)
    12 = @(5, 6)
    13 = @(6, 7)
    14 = @(7, 7)
    15 = @(7, 8)
    16 = @(7, 8)
    17 = @(6, 7)
    18 = @(7, 8)
    19 = @(6, 8)
    20 = @(5, 7)
    21 = @(4, 6)
    22 = @(8, 9)
    23 = @(6, 6)
    24 = @(2, 4)
    25 = @(8, 8)
    26 = @(7, 7)
    27 = @(7, 8)
    28 = @(7, 7)
    29 = @(8, 9)
    30 = @(7, 8)
    31 = @(7, 8)
    32 = @(7, 8)
    33 = @(7, 7)
    34 = @(5, 6)
    35 = @(1, 2)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
